$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row (becomes row 6) for the extra mtcars cylinder
#    bucket. This shifts the old rows 6,7,8,9 down to 7,8,9,10, and also
#    slides the existing merged cells (A6:A7 -> A7:A8, A9:G9 -> A10:G10)
#    and the sheet dimension (A1:H9 -> A1:H10) along with it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Insert()

# ---------------------------------------------------------------------------
# 2. The 4-cylinder bucket now spans two rows (5 and 6), so merge them and
#    top-align the cylinder count, matching the look already used for the
#    other two-row bucket (A7:A8).
# ---------------------------------------------------------------------------
$ws.Range("A5:A6").Merge()
$ws.Range("A5").VerticalAlignment = -4160
$ws.Range("A7").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 3. The new row 6 was cloned from row 5's formatting; repair the cells
#    whose formatting differs from row 5 (A6 blank "spacer" cell under the
#    merge, C6 "N" column, H6 trailing spacer) by pulling formats from cells
#    that already carry the right look elsewhere in the table. Use
#    Copy/PasteSpecial(formats) - a plain Copy(destination) onto a merged
#    cell silently breaks the merge, but PasteSpecial does not.
# ---------------------------------------------------------------------------
$ws.Range("H5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Row 5 (4-cylinder bucket) gets new summary numbers.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 91
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 2.14
$ws.Range("G5").ClearContents()

# ---------------------------------------------------------------------------
# 5. Row 6 (new automatic-transmission 4-cylinder bucket) gets its values.
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 81.8
$ws.Range("E6").Value = 21.87235698318771
$ws.Range("F6").Value = 2.3003
$ws.Range("G6").Value = 0.5982073312080948

# ---------------------------------------------------------------------------
# 6. Row 7 (6-cylinder bucket, manual) gets a refreshed N and Mean/SD
#    numbers.
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 131.6666666666667
$ws.Range("E7").Value = 37.52776749732568
$ws.Range("F7").Value = 2.755
$ws.Range("G7").Value = 0.1281600561797629

# ---------------------------------------------------------------------------
# 7. Row 8 (6-cylinder bucket, automatic) gets a new N and refreshed
#    Mean/SD numbers.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 115.25
$ws.Range("E8").Value = 9.178779875342908
$ws.Range("F8").Value = 3.38875
$ws.Range("G8").Value = 0.1162163929916946

# ---------------------------------------------------------------------------
# 8. Row 9 (8-cylinder bucket) gets a new N and refreshed Mean/SD numbers.
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 209.2142857142857
$ws.Range("E9").Value = 50.97688551827051
$ws.Range("F9").Value = 3.999214285714287
$ws.Range("G9").Value = 0.7594047444769265

# ---------------------------------------------------------------------------
# 9. Row 10 is the footer credit line, already carried down correctly by the
#    row insert/shift in step 1 - nothing further to do there.
# ---------------------------------------------------------------------------
